# CoordinateSystem.xlsx edit
# "Coord evaluation implemented + data format Cs, Point"
#
# The coordinate system table used to carry x / y / alpha plus co_x / co_y
# (the "is the coordinate inherited" flags). With coord evaluation now
# implemented, only a single numeric "x" value per coordinate system is
# kept (y and alpha are dropped), and co_x / co_y slide in right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "y" (E) and "alpha" (F) columns entirely; this shifts the
# former "co_x" (G) / "co_y" (H) columns left into E / F, right after "x".
$ws.Range("E1:F1").EntireColumn.Delete()

# CS_2's x value is recomputed under the new evaluation logic.
$ws.Range("D3").Value = -900
